# Update column G ("K" = strikeouts) values per row, replacing the old
# "Strike#" derived counts with the regenerated K counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4;  3  = 1;  4  = 1;  5  = 4;  6  = 1;  7  = 1;  8  = 0;  9  = 0;
    10 = 2;  11 = 3;  12 = 1;  13 = 1;  14 = 2;  15 = 2;  16 = 1;  17 = 2;
    18 = 2;  19 = 2;  20 = 0;  21 = 0;  22 = 2;  23 = 2;  24 = 4;  25 = 1;
    26 = 1;  27 = 1;  28 = 1;  29 = 1;  30 = 1;  31 = 1;  32 = 1;  33 = 0;
    34 = 1;  35 = 1;  36 = 0;  37 = 1;  38 = 1;  39 = 0;  40 = 0;  42 = 1;
    43 = 0;  44 = 1;  46 = 0;  47 = 0;  48 = 1;  49 = 0;  50 = 1;  51 = 1;
    52 = 2;  53 = 0;  54 = 1;  55 = 1;  56 = 2;  57 = 1;  58 = 1;  59 = 1;
    60 = 0;  61 = 0;  62 = 1;  63 = 0;  64 = 1;  65 = 0;  66 = 1;  68 = 1;
    69 = 1;  70 = 3;  71 = 1;  72 = 2;  74 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
